$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix block size error: B5 should be 15 instead of 30
$ws.Range("B5").Value = 15

# Force recalculation so dependent formulas (B8, B9, ...) update their cached values
$excel.Calculate()

# Update the selected cell / active cell as recorded in the saved view state
$ws.Range("E13").Select()
